$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.67"
$ws.Range("E2").Value = "'1.07%"
$ws.Range("D3").Value = "'41.46"
$ws.Range("E3").Value = "'4.54%"
$ws.Range("D4").Value = "'5.641"
$ws.Range("E4").Value = "'-1.11%"
$ws.Range("D5").Value = "'0.08313"
$ws.Range("E5").Value = "'3.52%"
$ws.Range("E6").Value = "'0.70%"
$ws.Range("D7").Value = "'8.769"
$ws.Range("E7").Value = "'1.46%"
$ws.Range("D8").Value = "'4.529"
$ws.Range("E8").Value = "'0.84%"
$ws.Range("D10").Value = "'0.9246"
$ws.Range("E10").Value = "'0.17%"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("E11").Value = "'1.73%"
$ws.Range("D12").Value = "'0.1959"
$ws.Range("E12").Value = "'-0.29%"
$ws.Range("D13").Value = "'0.09406"
$ws.Range("E13").Value = "'1.98%"
$ws.Range("D14").Value = "'0.03964"
$ws.Range("E14").Value = "'11.09%"
$ws.Range("E15").Value = "'1.09%"
$ws.Range("D16").Value = "'0.001305"
$ws.Range("E16").Value = "'1.72%"
$ws.Range("D17").Value = "'0.006149"
$ws.Range("E17").Value = "'-3.37%"
$ws.Range("D19").Value = "'3.444"
$ws.Range("D20").Value = "'0.3530"
$ws.Range("E20").Value = "'1.45%"
$ws.Range("D21").Value = "'8.387"
$ws.Range("E21").Value = "'-4.13%"
$ws.Range("E22").Value = "'1.67%"
$ws.Range("E23").Value = "'-1.42%"
$ws.Range("E24").Value = "'-0.40%"
$ws.Range("D25").Value = "'0.001257"
$ws.Range("E25").Value = "'-0.15%"
$ws.Range("D26").Value = "'0.004317"
$ws.Range("E26").Value = "'-6.51%"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'1.01%"
$ws.Range("D39").Value = "'0.02779"
$ws.Range("E39").Value = "'10.93%"
$ws.Range("D40").Value = "'0.05489"
$ws.Range("E40").Value = "'3.07%"
$ws.Range("D41").Value = "'0.007906"
$ws.Range("E41").Value = "'6.27%"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("E42").Value = "'1.15%"
$ws.Range("D43").Value = "'0.008932"
$ws.Range("E43").Value = "'-9.79%"
$ws.Range("D44").Value = "'0.002140"
$ws.Range("E44").Value = "'1.35%"
$ws.Range("D45").Value = "'0.01179"
$ws.Range("D46").Value = "'0.00007000"
$ws.Range("E46").Value = "'4.70%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.16%"
$ws.Range("D48").Value = "'0.003191"
$ws.Range("E48").Value = "'5.07%"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.16%"
